$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7818
$ws1.Range("F5").Value = 7818
$ws1.Range("F8").Value = 2122
$ws1.Range("F9").Value = 8597
$ws1.Range("F10").Value = 8597
$ws1.Range("F14").Value = 5744
$ws1.Range("F16").Value = 2718
$ws1.Range("F17").Value = 1196
$ws1.Range("F18").Value = 4606
$ws1.Range("F20").Value = 100
$ws1.Range("F22").Value = 601
$ws1.Range("F23").Value = 40
$ws1.Range("F24").Value = 3834
$ws1.Range("F25").Value = 77
$ws1.Range("F26").Value = 57
$ws1.Range("F27").Value = 48
$ws1.Range("F29").Value = 14
$ws1.Range("F30").Value = 5337
$ws1.Range("F31").Value = 4
$ws1.Range("F32").Value = 65
$ws1.Range("F34").Value = 383
$ws1.Range("F35").Value = 155
$ws1.Range("F36").Value = 378
$ws1.Range("F37").Value = 1755
$ws1.Range("F38").Value = 1004
$ws1.Range("F40").Value = 1117
$ws1.Range("F41").Value = 3735
$ws1.Range("F42").Value = 70
$ws1.Range("F44").Value = 32
$ws1.Range("F45").Value = 3454
$ws1.Range("F47").Value = 2330
$ws1.Range("F48").Value = 15
$ws1.Range("F50").Value = 471
$ws1.Range("F51").Value = 4

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 117
$ws2.Range("F3").Value = 142

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1358

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1358
$ws4.Range("F5").Value = 7818
$ws4.Range("F6").Value = 7818
$ws4.Range("F8").Value = 2122
$ws4.Range("F9").Value = 8597
$ws4.Range("F10").Value = 8597
$ws4.Range("F13").Value = 5744
$ws4.Range("F15").Value = 2718
$ws4.Range("F16").Value = 1196
$ws4.Range("F17").Value = 4606
$ws4.Range("F19").Value = 100
$ws4.Range("F20").Value = 117
$ws4.Range("F22").Value = 142
$ws4.Range("F23").Value = 601
$ws4.Range("F25").Value = 3834
$ws4.Range("F26").Value = 77
$ws4.Range("F27").Value = 57
$ws4.Range("F28").Value = 48
$ws4.Range("F30").Value = 14
$ws4.Range("F31").Value = 5337
$ws4.Range("F32").Value = 4
$ws4.Range("F33").Value = 65
$ws4.Range("F34").Value = 383
$ws4.Range("F35").Value = 155
$ws4.Range("F36").Value = 378
$ws4.Range("F38").Value = 1755
$ws4.Range("F39").Value = 1004
$ws4.Range("F41").Value = 1117
$ws4.Range("F43").Value = 3735
$ws4.Range("F44").Value = 70
$ws4.Range("F46").Value = 32
$ws4.Range("F47").Value = 3454
$ws4.Range("F48").Value = 2330
$ws4.Range("F49").Value = 471

Write-Output "Done applying updates"
